$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912" (sheet1): append new scrape rows 853-871, bump header info.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 31/12/2025 10:58:17"
$ws1.Range("A3").Value = "Total filas: 870"

$rows1 = @(
    @(853, "10:58:06", "11:01", "17_ROMERO", 3, "LP1912", "31/12/2025"),
    @(854, "10:58:06", "11:04", "23_HERNANDEZ", 6, "LP1912", "31/12/2025"),
    @(855, "10:58:06", "11:05", "14_ABASTO", 7, "LP1912", "31/12/2025"),
    @(856, "10:58:06", "11:09", "16_SANTA ANA", 11, "LP1912", "31/12/2025"),
    @(857, "10:58:06", "11:11", "15_ABASTO", 13, "LP1912", "31/12/2025"),
    @(858, "10:58:06", "11:14", "225_C ROCA-H SUR", 16, "LP1912", "31/12/2025"),
    @(859, "10:58:06", "11:25", "16_P MOR-SANTA ANA", 27, "LP1912", "31/12/2025"),
    @(860, "10:58:06", "11:26", "10_OLMOS", 28, "LP1912", "31/12/2025"),
    @(861, "10:58:06", "11:34", "10_OLMOS", 36, "LP1912", "31/12/2025"),
    @(862, "10:58:06", "11:34", "23_HERNANDEZ", 36, "LP1912", "31/12/2025"),
    @(863, "10:58:06", "11:40", "215A_EL PATO", 42, "LP1912", "31/12/2025"),
    @(864, "10:58:06", "11:45", "16_SANTA ANA", 47, "LP1912", "31/12/2025"),
    @(865, "10:58:06", "11:53", "15_ABASTO", 55, "LP1912", "31/12/2025"),
    @(866, "10:58:06", "11:54", "225_GOMEZ", 56, "LP1912", "31/12/2025"),
    @(867, "10:58:06", "11:57", "16_SANTA ANA", 59, "LP1912", "31/12/2025"),
    @(868, "10:58:06", "12:03", "23_HERNANDEZ", 65, "LP1912", "31/12/2025"),
    @(869, "10:58:06", "12:17", "15_ABASTO", 79, "LP1912", "31/12/2025"),
    @(870, "10:58:06", "12:18", "10_OLMOS", 80, "LP1912", "31/12/2025"),
    @(871, "10:58:06", "12:29", "215C_EL PATO", 91, "LP1912", "31/12/2025")
)

foreach ($r in $rows1) {
    $rowNum = $r[0]
    $ws1.Cells.Item($rowNum, 2).Value = $r[1]
    $ws1.Cells.Item($rowNum, 3).Value = $r[2]
    $ws1.Cells.Item($rowNum, 4).Value = $r[3]
    $ws1.Cells.Item($rowNum, 5).Value = $r[4]
    $ws1.Cells.Item($rowNum, 6).Value = $r[5]
    $ws1.Cells.Item($rowNum, 7).Value = $r[6]
}

# ---------------------------------------------------------------------------
# Sheet "LP1912-215" (sheet2): append new scrape rows 64-65, bump header info.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 31/12/2025 10:58:17"
$ws2.Range("A3").Value = "Total filas: 64"

$rows2 = @(
    @(64, "31/12/2025", "10:58:06", "11:40", "215A_EL PATO", 42, "LP1912"),
    @(65, "31/12/2025", "10:58:06", "12:29", "215C_EL PATO", 91, "LP1912")
)

foreach ($r in $rows2) {
    $rowNum = $r[0]
    $ws2.Cells.Item($rowNum, 2).Value = $r[1]
    $ws2.Cells.Item($rowNum, 3).Value = $r[2]
    $ws2.Cells.Item($rowNum, 4).Value = $r[3]
    $ws2.Cells.Item($rowNum, 5).Value = $r[4]
    $ws2.Cells.Item($rowNum, 6).Value = $r[5]
    $ws2.Cells.Item($rowNum, 7).Value = $r[6]
}

# ---------------------------------------------------------------------------
# Sheet "6203-6173" (sheet3): append new scrape row 103, bump header info.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 31/12/2025 10:58:17"
$ws3.Range("A3").Value = "Total filas: 102"

$ws3.Cells.Item(103, 2).Value = "31/12/2025"
$ws3.Cells.Item(103, 3).Value = "10:58:11"
$ws3.Cells.Item(103, 4).Value = "11:44"
$ws3.Cells.Item(103, 5).Value = "215C_LA PLATA"
$ws3.Cells.Item(103, 6).Value = 46
$ws3.Cells.Item(103, 7).Value = "L6203"
